$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map old emoji values to new replacement values (shared string table update).
# "+3" / "-3" are prefixed with a leading apostrophe so Excel stores them as
# literal text instead of silently coercing them to the numbers 3 / -3.
$map = @{
    "📘" = "⚠️"
    "📗" = "✅"
    "📙" = "'+3"
    "📕" = "'-3"
}

$used = $ws.UsedRange
$rows = $used.Rows.Count

for ($r = 1; $r -le $rows; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
